$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the FarmData2 School (FDS) assignment labels in column H:
# drop the space and leading zero from the assignment numbers (FDS 02 -> FDS2, etc.)
$ws.Range("H2").Value = "FDS2, due 1/31"
$ws.Range("H3").Value = "FDS3, due 2/7"
$ws.Range("H4").Value = "FDS4, due 2/14"
$ws.Range("H5").Value = "FDS5, due 2/21"
$ws.Range("H6").Value = "FDS6, due 2/28"
$ws.Range("H7").Value = "FDS7, due 3/7"

# Update the saved view state: scroll the frozen pane down a bit and move
# the active selection (matches the author's new cursor position).
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 5
$ws.Range("F19").Select() | Out-Null
